$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
$ws.Activate()

$ws.Range("G3").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0876_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nSwitchApp(NATIVE_APP);`nDrawSignature(signature_view_xpath);`nTakeScreenshot(VT200-0876);`nClickNativeIcon(signature_ok_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nwait(2);`nvalidate4;"
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0877_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nSwitchApp(NATIVE_APP);`nTakeScreenshot(VT200-0877);`nClickNativeIcon(signature_cancel_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nwait(2);`nvalidate4;"
$ws.Range("G7").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0880_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nSwitchApp(NATIVE_APP);`nDrawSignature(signature_view_xpath);`nTakeScreenshot(VT200-0880);`nClickNativeIcon(signature_ok_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nwait(2);`nvalidate4;"
$ws.Range("G8").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0881_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nSwitchApp(NATIVE_APP);`nDrawSignature(signature_view_xpath);`nTakeScreenshot(VT200-0881-01);`nClickNativeIcon(signature_ok_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nwait(2);`nTakeScreenshot(VT200-0881-02);`nvalidate4;"
$ws.Range("G12").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0885_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nRotate_Screen(landscape);`nwait(2);`nTakeScreenshot(VT200-0885);`nwait(2);`nvalidate4;"
$ws.Range("G13").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0886_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nRotate_Screen(landscape);`nwait(2);`nTakeScreenshot(VT200-0886);`nwait(2);`nvalidate4;"
$ws.Range("G14").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0887_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nSwitchApp(NATIVE_APP);`nDrawSignature(signature_view_xpath);`nTakeScreenshot(VT200-0887-01);`nClickNativeIcon(signature_clear_xpath);`nwait(2);`nTakeScreenshot(VT200-0887-02);`nClickNativeIcon(signature_ok_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nwait(2);`nvalidate4;"
$ws.Range("G16").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0889_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nTakeScreenshot(VT200-0889-01);`nLock_UnlockScreen(lock);`nLock_UnlockScreen(unlock);`nwait(2);`nTakeScreenshot(VT200-0889-02);`nvalidate4;"
$ws.Range("G20").Value = "wait(3);`nvalidate1;`nlink_Click(signature_test_link);`nvalidate2;`nSelectTestToRun(VT200_0893_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nRotate_Screen(landscape);`nwait(2);`nTakeScreenshot(VT200-0893-01);`nwait(2);`nRotate_Screen(potrait);`nwait(2);`nTakeScreenshot(VT200-0893-02);`nwait(2);`nvalidate4;"

$ws.Range("G3").Select()
